$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new tracking number, new actual rate -> now a FAIL (mismatch vs expected rate)
$ws.Range("P2").Value = "'" + "310109788588"
$ws.Range("P2").Style = "Normal"
$ws.Range("Q2").Value = "'" + '$19.13'
$ws.Range("Q2").Style = "Normal"
$ws.Range("R2").Value = "FAIL"

# Row 3: new tracking number, new actual rate -> now a FAIL (mismatch vs expected rate)
$ws.Range("P3").Value = "'" + "310109788599"
$ws.Range("P3").Style = "Normal"
$ws.Range("Q3").Value = "'" + '$27.63'
$ws.Range("Q3").Style = "Normal"
$ws.Range("R3").Value = "FAIL"

# Rows 16-21: results corrected from FAIL to PASS
$ws.Range("R16").Value = "PASS"
$ws.Range("R17").Value = "PASS"
$ws.Range("R18").Value = "PASS"
$ws.Range("R19").Value = "PASS"
$ws.Range("R20").Value = "PASS"
$ws.Range("R21").Value = "PASS"
